$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11 (write A11 first so its shared string is interned before D1's)
$ws.Range("A11").Value = "Walfaanaa Magarsaa"
$ws.Range("B11").Value = 912861288
$ws.Range("C11").Value = 10000
$ws.Range("D11").Value = 100000

# New column D header + data
$ws.Range("D1").Value = "Total money for winners"
$ws.Range("D2:D10").Value = 100000

# Formatting for D1 header: inherit the A1:C1 header look (bold font,
# centered/top alignment, thin border) then drop the top/bottom edges so
# only the left/right border remains.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Borders.Item(8).LineStyle = -4142
$ws.Range("D1").Borders.Item(9).LineStyle = -4142
$ws.Application.CutCopyMode = $false

# Column widths (autofit to match bestFit behaviour of columns A/B)
$ws.Columns.Item(3).EntireColumn.AutoFit()
$ws.Columns.Item(4).EntireColumn.AutoFit()

$ws.Range("E8").Select()
